$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (D = Price, E = Volume) keep their
# original Text cell type instead of Excel auto-converting numeric-
# looking strings (e.g. "0.740", "75.30") into Number cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.358.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.728"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "256.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.20%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.67"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.37"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0767"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0989"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.31"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.196.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.740"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.909.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.346.61"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.30"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.50%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.15"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.00%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.86"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.67%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.39"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.67"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +25.12%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +14.95%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.12%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.924"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "100.62"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +12.11%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0222"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.55%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.21%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0650"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.349.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.21%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.77"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.77"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.47"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.13%  "
